# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp cell.
# - Bielorrusia overtakes Honduras in total cases, so the two countries swap
#   their row position in the (already case-sorted) table.
# - Refresh the day's case/death counters for the affected countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 15:03"

# Bielorrusia / Honduras swap places (row 54 <-> row 55)
$ws.Range("A54").Value = "Bielorrusia"
$ws.Range("A55").Value = "Honduras"

# Per-row numeric updates: row number -> @{ col = value }
$updates = @{
    5   = @{ B = 7244024; C = 6942; D = 6304748; E = 828590;            G = 69; H = 110686 }
    23  = @{ B = 340590;  C = 501;  D = 326820;  E = 8662;               G = 21; H = 5108 }
    29  = @{ B = 196163;  C = 7287;                                      G = 32; H = 6663 }
    47  = @{ B = 101332;                                                 G = 6;  H = 5907 }
    54  = @{ B = 85121;   C = 597;  D = 78218;   E = 5992;                G = 5;  H = 911 }
    55  = @{ B = 84852;   C = 439;  D = 32772;   E = 49552;               G = 7;  H = 2528 }
    73  = @{ B = 42750;   C = 369;  D = 39570;   E = 2564;                G = 4;  H = 616 }
    76  = @{                                     E = 5159;               G = 1;  H = 1481 }
    79  = @{ B = 33593;   C = 492;  D = 27680;   E = 5238;                G = 1;  H = 675 }
    81  = @{ B = 30766;   C = 286;  D = 26087;   E = 3775;                G = 5;  H = 904 }
    133 = @{ B = 5055;    C = 17;   D = 3357;    E = 1685 }
    154 = @{ B = 2478;    C = 67;   D = 2330;    E = 107 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

Write-Output "done"
